# Updated symbol list on Sun Jan 15 13:10:39 UTC 2023 with GitHub Actions
# Refresh Price / Volume(1h) / Hora columns for each coin row.
# Values are entered with a leading apostrophe so Excel stores them as
# literal text (matching the existing text-formatted cells) instead of
# coercing number/percent-looking strings into numeric values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'297.57"
$ws.Range("E2").Value = "'-2.29%"
$ws.Range("G2").Value = "'13"
$ws.Range("D3").Value = "'31.21"
$ws.Range("E3").Value = "'-2.44%"
$ws.Range("G3").Value = "'13"
$ws.Range("D4").Value = "'5.125"
$ws.Range("E4").Value = "'-2.51%"
$ws.Range("G4").Value = "'13"
$ws.Range("D5").Value = "'0.07321"
$ws.Range("E5").Value = "'-3.28%"
$ws.Range("G5").Value = "'13"
$ws.Range("D6").Value = "'7.747"
$ws.Range("E6").Value = "'-1.31%"
$ws.Range("G6").Value = "'13"
$ws.Range("D7").Value = "'1.746"
$ws.Range("E7").Value = "'17.62%"
$ws.Range("G7").Value = "'13"
$ws.Range("D8").Value = "'3.726"
$ws.Range("E8").Value = "'-0.82%"
$ws.Range("G8").Value = "'13"
$ws.Range("D9").Value = "'0.9245"
$ws.Range("E9").Value = "'1.09%"
$ws.Range("G9").Value = "'13"
$ws.Range("D10").Value = "'0.1669"
$ws.Range("E10").Value = "'-1.46%"
$ws.Range("G10").Value = "'13"
$ws.Range("D11").Value = "'0.06942"
$ws.Range("E11").Value = "'-7.64%"
$ws.Range("G11").Value = "'13"
$ws.Range("D12").Value = "'0.07979"
$ws.Range("E12").Value = "'-1.06%"
$ws.Range("G12").Value = "'13"
$ws.Range("D13").Value = "'0.03009"
$ws.Range("E13").Value = "'0.45%"
$ws.Range("G13").Value = "'13"
$ws.Range("E14").Value = "'0.17%"
$ws.Range("G14").Value = "'13"
$ws.Range("D15").Value = "'0.001512"
$ws.Range("E15").Value = "'0.30%"
$ws.Range("G15").Value = "'13"
$ws.Range("D16").Value = "'0.006117"
$ws.Range("E16").Value = "'-1.85%"
$ws.Range("G16").Value = "'13"
$ws.Range("E17").Value = "'-1.12%"
$ws.Range("G17").Value = "'13"
$ws.Range("D18").Value = "'2.218"
$ws.Range("E18").Value = "'-0.49%"
$ws.Range("G18").Value = "'13"
$ws.Range("D19").Value = "'0.3226"
$ws.Range("E19").Value = "'-2.58%"
$ws.Range("G19").Value = "'13"
$ws.Range("D20").Value = "'0.1347"
$ws.Range("E20").Value = "'0.32%"
$ws.Range("G20").Value = "'13"
$ws.Range("D21").Value = "'4.585"
$ws.Range("E21").Value = "'2.32%"
$ws.Range("G21").Value = "'13"
$ws.Range("D22").Value = "'0.04651"
$ws.Range("E22").Value = "'2.44%"
$ws.Range("G22").Value = "'13"
$ws.Range("D23").Value = "'0.1582"
$ws.Range("E23").Value = "'-2.77%"
$ws.Range("G23").Value = "'13"
$ws.Range("E24").Value = "'0.68%"
$ws.Range("G24").Value = "'13"
$ws.Range("D25").Value = "'0.004743"
$ws.Range("E25").Value = "'6.70%"
$ws.Range("G25").Value = "'13"
$ws.Range("E26").Value = "'-17.98%"
$ws.Range("G26").Value = "'13"
$ws.Range("D27").Value = "'0.0001873"
$ws.Range("E27").Value = "'7.85%"
$ws.Range("G27").Value = "'13"
$ws.Range("G28").Value = "'13"
$ws.Range("G29").Value = "'13"
$ws.Range("G30").Value = "'13"
$ws.Range("G31").Value = "'13"
$ws.Range("G32").Value = "'13"
$ws.Range("G33").Value = "'13"
$ws.Range("G34").Value = "'13"
$ws.Range("G35").Value = "'13"
$ws.Range("G36").Value = "'13"
$ws.Range("G37").Value = "'13"
$ws.Range("G38").Value = "'13"
$ws.Range("D39").Value = "'0.01707"
$ws.Range("E39").Value = "'2.78%"
$ws.Range("G39").Value = "'13"
$ws.Range("D40").Value = "'0.04429"
$ws.Range("E40").Value = "'-1.97%"
$ws.Range("G40").Value = "'13"
$ws.Range("D41").Value = "'0.007137"
$ws.Range("E41").Value = "'-0.97%"
$ws.Range("G41").Value = "'13"
$ws.Range("D42").Value = "'0.1331"
$ws.Range("E42").Value = "'-1.35%"
$ws.Range("G42").Value = "'13"
$ws.Range("D43").Value = "'0.002187"
$ws.Range("E43").Value = "'-2.36%"
$ws.Range("G43").Value = "'13"
$ws.Range("D44").Value = "'0.01110"
$ws.Range("E44").Value = "'-14.93%"
$ws.Range("G44").Value = "'13"
$ws.Range("D45").Value = "'0.00006084"
$ws.Range("E45").Value = "'-1.80%"
$ws.Range("G45").Value = "'13"
$ws.Range("E46").Value = "'-21.24%"
$ws.Range("G46").Value = "'13"
$ws.Range("E47").Value = "'170.81%"
$ws.Range("G47").Value = "'13"
$ws.Range("G48").Value = "'13"
$ws.Range("G49").Value = "'13"
$ws.Range("G50").Value = "'13"
$ws.Range("G51").Value = "'13"
